# Auto-generated: bulk market-price / profit recompute refresh
# (mirrors the scheduled-runner "chore: update Sheets" commit)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5436384
$ws.Range("I19").Value = 8929388
$ws.Range("J19").Value = 2822.889
$ws.Range("K19").Value = 8929388
$ws.Range("L19").Value = 2822.889
$ws.Range("M19").Value = -8929213
$ws.Range("N19").Value = -3172.889
$ws.Range("H62").Value = 48618520
$ws.Range("I62").Value = 19238846
$ws.Range("J62").Value = 125005670
$ws.Range("K62").Value = 19238846
$ws.Range("L62").Value = 125005670
$ws.Range("M62").Value = -19238222
$ws.Range("N62").Value = -125006918
$ws.Range("H65").Value = 48618520
$ws.Range("I65").Value = 19238846
$ws.Range("J65").Value = 125005670
$ws.Range("K65").Value = 96194230
$ws.Range("L65").Value = 625028350
$ws.Range("M65").Value = -96191110
$ws.Range("N65").Value = -625034590
$ws.Range("H132").Value = 1523172.1
$ws.Range("I132").Value = 822.4091
$ws.Range("J132").Value = 15876756
$ws.Range("K132").Value = 2467.2273
$ws.Range("L132").Value = 47630268
$ws.Range("M132").Value = 62.77269999999999
$ws.Range("N132").Value = -47635328
$ws.Range("H135").Value = 690.0833
$ws.Range("I135").Value = 563
$ws.Range("J135").Value = 1173
$ws.Range("K135").Value = 5067
$ws.Range("L135").Value = 10557
$ws.Range("M135").Value = -2532
$ws.Range("N135").Value = -15627

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1097781
$ws.Range("I61").Value = 1437588.5
$ws.Range("J61").Value = 2845.7778
$ws.Range("K61").Value = 1437588.5
$ws.Range("L61").Value = 2845.7778
$ws.Range("M61").Value = -1437376.5
$ws.Range("N61").Value = -3269.7778
$ws.Range("H63").Value = 2069.2307
$ws.Range("I63").Value = 1908.3334
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 1908.3334
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1222.3334
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 2069.2307
$ws.Range("I66").Value = 1908.3334
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 9541.666999999999
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -6109.666999999999
$ws.Range("N66").Value = -26864
$ws.Range("H74").Value = 69699110
$ws.Range("I74").Value = 62500980
$ws.Range("J74").Value = 88894130
$ws.Range("K74").Value = 62500980
$ws.Range("L74").Value = 88894130
$ws.Range("M74").Value = -62500106
$ws.Range("N74").Value = -88895878
$ws.Range("H77").Value = 69699110
$ws.Range("I77").Value = 62500980
$ws.Range("J77").Value = 88894130
$ws.Range("K77").Value = 312504900
$ws.Range("L77").Value = 444470650
$ws.Range("M77").Value = -312500532
$ws.Range("N77").Value = -444479386
$ws.Range("H132").Value = 10103641
$ws.Range("I132").Value = 11366535
$ws.Range("J132").Value = 5052064
$ws.Range("K132").Value = 34099605
$ws.Range("L132").Value = 15156192
$ws.Range("M132").Value = -34097075
$ws.Range("N132").Value = -15161252
$ws.Range("H136").Value = 1097781
$ws.Range("I136").Value = 1437588.5
$ws.Range("J136").Value = 2845.7778
$ws.Range("K136").Value = 4312765.5
$ws.Range("L136").Value = 8537.3334
$ws.Range("M136").Value = -4310215.5
$ws.Range("N136").Value = -13637.3334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17867294
$ws.Range("I20").Value = 27784582
$ws.Range("J20").Value = 16174.8
$ws.Range("K20").Value = 27784582
$ws.Range("L20").Value = 16174.8
$ws.Range("M20").Value = -27784335
$ws.Range("N20").Value = -16668.8
$ws.Range("H94").Value = 1262.7391
$ws.Range("I94").Value = 870.6842
$ws.Range("J94").Value = 3125
$ws.Range("K94").Value = 870.6842
$ws.Range("L94").Value = 3125
$ws.Range("M94").Value = -419.6842
$ws.Range("N94").Value = -4027
$ws.Range("H134").Value = 10303225
$ws.Range("I134").Value = 14286672
$ws.Range("J134").Value = 2102011.5
$ws.Range("K134").Value = 42860016
$ws.Range("L134").Value = 6306034.5
$ws.Range("M134").Value = -42857481
$ws.Range("N134").Value = -6311104.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1589.9
$ws.Range("I16").Value = 1666.5
$ws.Range("J16").Value = 1475
$ws.Range("K16").Value = 1666.5
$ws.Range("L16").Value = 1475
$ws.Range("M16").Value = -1379.5
$ws.Range("N16").Value = -2049
$ws.Range("H107").Value = 536.76086
$ws.Range("I107").Value = 208.64706
$ws.Range("J107").Value = 729.10345
$ws.Range("K107").Value = 208.64706
$ws.Range("L107").Value = 729.10345
$ws.Range("M107").Value = 1711.35294
$ws.Range("N107").Value = -4569.10345
$ws.Range("H113").Value = 1589.9
$ws.Range("I113").Value = 1666.5
$ws.Range("J113").Value = 1475
$ws.Range("K113").Value = 1666.5
$ws.Range("L113").Value = 1475
$ws.Range("M113").Value = 503.5
$ws.Range("N113").Value = -5815
$ws.Range("H134").Value = 680880
$ws.Range("I134").Value = 2507.2678
$ws.Range("K134").Value = 7521.803400000001
$ws.Range("M134").Value = -4986.803400000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 529
$ws.Range("J68").Value = 660.3333
$ws.Range("L68").Value = 1980.9999
$ws.Range("N68").Value = -3602.9999
$ws.Range("H71").Value = 529
$ws.Range("J71").Value = 660.3333
$ws.Range("L71").Value = 5942.9997
$ws.Range("N71").Value = -14054.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9263559
$ws.Range("I70").Value = 3476747
$ws.Range("J70").Value = 55558056
$ws.Range("K70").Value = 3476747
$ws.Range("L70").Value = 55558056
$ws.Range("M70").Value = -3476477
$ws.Range("N70").Value = -55558596
$ws.Range("H73").Value = 9263559
$ws.Range("I73").Value = 3476747
$ws.Range("J73").Value = 55558056
$ws.Range("K73").Value = 3476747
$ws.Range("L73").Value = 55558056
$ws.Range("M73").Value = -3475811
$ws.Range("N73").Value = -55559928
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H141").Value = 49500
$ws.Range("J141").Value = 49500
$ws.Range("L141").Value = 49500
$ws.Range("N141").Value = -59860

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 15734.538
$ws.Range("I107").Value = 20165
$ws.Range("J107").Value = 966.3333
$ws.Range("K107").Value = 60495
$ws.Range("L107").Value = 2898.9999
$ws.Range("M107").Value = -58575
$ws.Range("N107").Value = -6738.9999
$ws.Range("H136").Value = 2833.2632
$ws.Range("I136").Value = 681.5227
$ws.Range("J136").Value = 10116.077
$ws.Range("K136").Value = 2044.5681
$ws.Range("L136").Value = 30348.231
$ws.Range("M136").Value = 505.4319
$ws.Range("N136").Value = -35448.231
